$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '96.071.77'
$ws.Range("E2").Value = '  -0.88%  '

Set-TextValue $ws.Range("D3") '3.624.84'
$ws.Range("E3").Value = '  -1.63%  '

Set-TextValue $ws.Range("D4") '2.75'
$ws.Range("E4").Value = '  +30.91%  '

$ws.Range("E5").Value = '  +0.11%  '

Set-TextValue $ws.Range("D6") '224.63'
$ws.Range("E6").Value = '  -4.49%  '

$ws.Range("E7").Value = '  -1.91%  '

Set-TextValue $ws.Range("D8") '0.425'
$ws.Range("E8").Value = '  -1.51%  '

Set-TextValue $ws.Range("D9") '1.21'
$ws.Range("E9").Value = '  +11.97%  '

$ws.Range("E10").Value = '  +0.06%  '

Set-TextValue $ws.Range("D11") '3.620.30'
$ws.Range("E11").Value = '  -1.72%  '

Set-TextValue $ws.Range("D12") '48.49'
$ws.Range("E12").Value = '  +9.51%  '

$ws.Range("E13").Value = '  +4.94%  '

$ws.Range("E14").Value = '  -5.34%  '

Set-TextValue $ws.Range("D15") '6.56'
$ws.Range("E15").Value = '  -3.29%  '

Set-TextValue $ws.Range("D16") '4.302.35'
$ws.Range("E16").Value = '  -1.63%  '

Set-TextValue $ws.Range("D17") '95.786.95'
$ws.Range("E17").Value = '  -0.85%  '

Set-TextValue $ws.Range("D18") '24.54'
$ws.Range("E18").Value = '  +32.63%  '

$ws.Range("E19").Value = '  +2.73%  '

Set-TextValue $ws.Range("D20") '14.00'
$ws.Range("E20").Value = '  +8.62%  '

Set-TextValue $ws.Range("D21") '3.629.94'
$ws.Range("E21").Value = '  -1.69%  '

Set-TextValue $ws.Range("D22") '0.295'
$ws.Range("E22").Value = '  +46.40%  '

Set-TextValue $ws.Range("D23") '0.539'
$ws.Range("E23").Value = '  +1.34%  '

Set-TextValue $ws.Range("D24") '135.47'
$ws.Range("E24").Value = '  +23.95%  '

Set-TextValue $ws.Range("D25") '524.74'
$ws.Range("E25").Value = '  +1.88%  '

$ws.Range("E26").Value = '  -4.20%  '

Set-TextValue $ws.Range("D27") '0.0000203'
$ws.Range("E27").Value = '  -7.41%  '

Set-TextValue $ws.Range("D28") '6.86'
$ws.Range("E28").Value = '  +0.13%  '

Set-TextValue $ws.Range("D29") '3.805.85'
$ws.Range("E29").Value = '  -2.05%  '

Set-TextValue $ws.Range("D30") '12.92'
$ws.Range("E30").Value = '  -2.85%  '

Set-TextValue $ws.Range("D31") '13.30'
$ws.Range("E31").Value = '  +6.62%  '

Set-TextValue $ws.Range("D32") '3.14'
$ws.Range("E32").Value = '  +5.32%  '

$ws.Range("E33").Value = '  +0.15%  '

Set-TextValue $ws.Range("D34") '0.634'
$ws.Range("E34").Value = '  +8.05%  '

$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range("D35") '33.36'
$ws.Range("E35").Value = '  +2.93%  '

$ws.Range("B36").Value = 'Cronos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range("D36") '0.184'
$ws.Range("E36").Value = '  -1.55%  '

$ws.Range("B37").Value = 'Fetch.AI'
$ws.Range("C37").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range("D37") '1.81'
$ws.Range("E37").Value = '  +0.54%  '

$ws.Range("B38").Value = 'Binance-PegBSC-USD'
$ws.Range("C38").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue $ws.Range("D38") '1.00'
$ws.Range("E38").Value = '  +0.33%  '

$ws.Range("B39").Value = 'Algorand'
$ws.Range("C39").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range("D39") '0.540'
$ws.Range("E39").Value = '  +10.41%  '

$ws.Range("B40").Value = 'USDe'
$ws.Range("C40").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws.Range("D40") '1.00'
$ws.Range("E40").Value = '  -0.03%  '

$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D41") '7.30'
$ws.Range("E41").Value = '  +8.57%  '

Set-TextValue $ws.Range("D42") '592.34'
$ws.Range("E42").Value = '  -5.67%  '

Set-TextValue $ws.Range("D43") '8.37'
$ws.Range("E43").Value = '  -3.38%  '

Set-TextValue $ws.Range("D44") '0.0535'
$ws.Range("E44").Value = '  +19.99%  '

Set-TextValue $ws.Range("D45") '41.67'
$ws.Range("E45").Value = '  +4.53%  '

Set-TextValue $ws.Range("D46") '1.01'
$ws.Range("E46").Value = '  +6.57%  '

$ws.Range("E47").Value = '  -4.92%  '

Set-TextValue $ws.Range("D48") '1.98'
$ws.Range("E48").Value = '  -1.46%  '

$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range("D49") '9.25'
$ws.Range("E49").Value = '  +7.34%  '

$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range("D50") '238.21'
$ws.Range("E50").Value = '  +17.18%  '

Set-TextValue $ws.Range("D51") '2.31'
$ws.Range("E51").Value = '  -2.35%  '
